$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$updates = @(
    @{ Row = 2; D = "29.314.76"; E = "  +0.56%  " }
    @{ Row = 3; D = "1.873.33"; E = "  +0.62%  " }
    @{ Row = 5; D = "0.7114"; E = "  +0.35%  " }
    @{ Row = 6; D = $null; E = "  +0.16%  " }
    @{ Row = 7; D = $null; E = "  +0.07%  " }
    @{ Row = 8; D = "0.3103"; E = "  +0.02%  " }
    @{ Row = 9; D = "0.07759"; E = "  +1.61%  " }
    @{ Row = 10; D = "25.02"; E = "  +1.40%  " }
    @{ Row = 11; D = "0.08396"; E = "  +0.48%  " }
    @{ Row = 12; D = "1.859.99"; E = "  +0.20%  " }
    @{ Row = 13; D = "5.236"; E = "  +0.82%  " }
    @{ Row = 14; D = "0.7110"; E = "  +0.51%  " }
    @{ Row = 15; D = $null; E = "  -0.02%  " }
    @{ Row = 16; D = "29.320.39"; E = "  +0.53%  " }
    @{ Row = 17; D = $null; E = "  +2.59%  " }
    @{ Row = 18; D = "0.000008178"; E = "  +4.89%  " }
    @{ Row = 19; D = "239.33"; E = "  -1.39%  " }
    @{ Row = 20; D = $null; E = "  +0.80%  " }
    @{ Row = 21; D = "2.119.97"; E = "  +0.56%  " }
    @{ Row = 22; D = "1.001"; E = "  +0.12%  " }
    @{ Row = 23; D = "7.752"; E = "  -1.36%  " }
    @{ Row = 24; D = $null; E = "  +0.12%  " }
    @{ Row = 25; D = "0.1594"; E = "  +0.71%  " }
    @{ Row = 26; D = "162.66"; E = "  -0.67%  " }
    @{ Row = 27; D = "9.018"; E = "  +0.67%  " }
    @{ Row = 28; D = "18.50"; E = "  +0.63%  " }
    @{ Row = 29; D = $null; E = "  +0.63%  " }
    @{ Row = 30; D = "4.404"; E = "  +0.20%  " }
    @{ Row = 31; D = $null; E = "  +1.57%  " }
    @{ Row = 32; D = "1.287"; E = "  -2.56%  " }
    @{ Row = 33; D = "0.05289"; E = "  +2.62%  " }
    @{ Row = 34; D = $null; E = "  +1.33%  " }
    @{ Row = 35; D = "1.176"; E = "  +1.07%  " }
    @{ Row = 36; D = "0.7448"; E = "  -7.00%  " }
    @{ Row = 37; D = $null; E = "  +0.60%  " }
    @{ Row = 38; D = "0.01876"; E = "  +1.71%  " }
    @{ Row = 39; D = "1.220.60"; E = "  +4.78%  " }
    @{ Row = 40; D = "2.724"; E = "  +1.18%  " }
    @{ Row = 41; D = "6.513"; E = "  +4.90%  " }
    @{ Row = 42; D = "109.66"; E = "  +7.38%  " }
    @{ Row = 43; D = "0.8861"; E = "  -0.34%  " }
    @{ Row = 44; D = "72.33"; E = "  -0.62%  " }
    @{ Row = 46; D = "2.018.50"; E = "  +0.56%  " }
    @{ Row = 47; D = "1.796"; E = "  +1.11%  " }
    @{ Row = 48; D = "0.5195"; E = "  -0.14%  " }
    @{ Row = 49; D = "0.00000000122"; E = "  +2.34%  " }
    @{ Row = 50; D = "9.351"; E = "  +0.23%  " }
    @{ Row = 51; D = $null; E = "  +1.00%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextCell "D$($u.Row)" $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
